# Append newly-predicted weather readings (Time Stamp, Temperature,
# Humidity, Pressure, Wind Speed) below the existing rows on the
# "Current Weather" sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("2017.05.20 13.39.32", 304.8999938964844,  25.0, 1011.0, 5.099999904632568),
    @("2017.05.20 14.11.10", 304.8900146484375,  25.0, 1011.0, 5.099999904632568),
    @("2017.05.20 14.14.51", 31.15999984741211,  26.0, 1011.0, 4.599999904632568),
    @("2017.05.21 11.49.32", 27.0,               57.0, 1013.0, 3.0999999046325684),
    @("2017.05.21 12.18.31", 27.829999923706055, 47.0, 1013.0, 4.599999904632568),
    @("2017.05.21 12.19.17", 27.829999923706055, 47.0, 1013.0, 4.599999904632568),
    @("2017.05.21 12.22.03", 27.829999923706055, 47.0, 1013.0, 4.599999904632568),
    @("2017.05.21 12.29.27", 27.829999923706055, 47.0, 1013.0, 4.599999904632568),
    @("2017.05.21 12.32.26", 27.829999923706055, 47.0, 1013.0, 4.599999904632568),
    @("2017.05.21 12.34.57", 27.829999923706055, 47.0, 1013.0, 4.599999904632568),
    @("2017.05.21 12.45.55", 27.829999923706055, 47.0, 1013.0, 4.599999904632568)
)

# Existing data occupies rows 1 (header) and 2-3; new rows continue from row 4.
$startRow = 4

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = $startRow + $i
    $data = $newRows[$i]

    $ws.Cells.Item($row, 1).Value = $data[0]
    $ws.Cells.Item($row, 2).Value = $data[1]
    $ws.Cells.Item($row, 3).Value = $data[2]
    $ws.Cells.Item($row, 4).Value = $data[3]
    $ws.Cells.Item($row, 5).Value = $data[4]
}
